$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9927031993865967
$ws.Range("B1").Value = 1.610569953918457
$ws.Range("C1").Value = 2.739323854446411
$ws.Range("D1").Value = 1.537822604179382
$ws.Range("E1").Value = 0.8235983848571777
